# Apply trade #6 close to the "live_trading_results" workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.3     # Total P&L %
$wsSummary.Range("B6").Value = 6       # Total Trades
$wsSummary.Range("B9").Value = 50      # Win Rate %

# --- Sheet "Strategy Status" -------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 6        # MarketMaking Trades
$wsStatus.Range("G4").Value = 50       # MarketMaking Win Rate %

# --- Helper: append the closed trade #6 row to a trades sheet ----------
function Add-TradeSixRow($ws) {
    $ws.Range("A7").Value = 6

    # The Date/Time columns hold plain text that looks like a date/time
    # (e.g. "2026-02-17"). Setting .Value directly would make Excel
    # auto-convert it to a date/time serial number, which does not match
    # how the rest of the sheet stores these columns (as text). Force the
    # cell to text ("@") before assigning, then restore the style back to
    # Normal so the cell ends up identical (plain text, default style) to
    # its siblings in the column.
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = "2026-02-17"
    $ws.Range("B7").Style = "Normal"

    $ws.Range("C7").NumberFormat = "@"
    $ws.Range("C7").Value = "19:43:08"
    $ws.Range("C7").Style = "Normal"
    $ws.Range("D7").Value = "MarketMaking"
    $ws.Range("E7").Value = "DOWN"
    $ws.Range("F7").Value = 0.39
    $ws.Range("G7").Value = 0.39
    $ws.Range("H7").Value = "CLOSED"
    $ws.Range("I7").Value = 0
    $ws.Range("J7").Value = 0
    $ws.Range("K7").Value = 100.09
    $ws.Range("L7").Value = 0
    $ws.Range("M7").Value = 0
    $ws.Range("N7").Value = 0.6
    $ws.Range("O7").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P7").Value = "early_exit"
    $ws.Range("Q7").Value = 0.13
}

# --- Sheet "All Trades" -------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeSixRow $wsAllTrades

# --- Sheet "MarketMaking" ------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeSixRow $wsMarketMaking
